# Refresh the cryptos list: update each coin's Price (column D) and
# Volume(1h) (column E) cell to the latest scraped figures.
#
# Several Price values are plain decimal numbers (e.g. "581.85"); setting
# .Value on them would let Excel auto-convert the cell to a real number
# (losing the text formatting the sheet relies on), so for those cells we
# briefly force a Text number format, assign the value, then clear the
# format again so the cell ends up as plain text with the default style
# (matching the rest of the sheet, which never set an explicit NumberFormat).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.894.70'
$ws.Range("D3").Value = '3.282.15'
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.85'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.48%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("E9").Value = '  +7.63%  '
$ws.Range("E10").Value = '  +1.76%  '
$ws.Range("E11").Value = '  +6.22%  '
$ws.Range("D12").Value = '3.853.33'
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.81'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.86%  '
$ws.Range("D15").Value = '67.860.95'
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("E16").Value = '  +3.71%  '
$ws.Range("D17").Value = '3.282.01'
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("E18").Value = '  +2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.55'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '377.58'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.69'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.61%  '
$ws.Range("E24").Value = '  +3.87%  '
$ws.Range("E25").Value = '  +5.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.80'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.182'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.93%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  +3.21%  '
$ws.Range("E30").Value = '  +6.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.98'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.28'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +7.16%  '
$ws.Range("E34").Value = '  +5.34%  '
$ws.Range("E35").Value = '  +5.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.855'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.18%  '
$ws.Range("E38").Value = '  +2.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.10'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.86%  '
$ws.Range("E40").Value = '  +10.13%  '
$ws.Range("E41").Value = '  +10.52%  '
$ws.Range("E42").Value = '  +4.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.87'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +8.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '352.01'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +7.33%  '
$ws.Range("D45").Value = '2.665.18'
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0683'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.47%  '
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("E49").Value = '  +5.46%  '
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.17'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.68%  '
